<#
.SYNOPSIS
    Generalizes the "Hämtlista" pickup-list template by removing the
    hard-coded, one-off event title ("Påskmarknad 2022") that used to
    live in cell B2, and resets the sheet's active selection to B2 so the
    template opens in a clean, predictable state for the next event.

.DESCRIPTION
    The template header cell (B2) previously contained a literal string
    referencing a specific market ("Påskmarknad 2022"). Since the
    template is reused for many different events, that text is removed
    here (the cell keeps its existing formatting/style, only its
    content is cleared) so the template is no longer tied to a single
    occasion. The now-unused shared string is dropped automatically by
    the workbook when no cell references it any more.

    Basic error handling is included so that a failure on one sheet does
    not leave the workbook in a half-edited state without at least a
    clear diagnostic message.
#>

try {
    $wb = $excel.ActiveWorkbook
    if (-not $wb) {
        throw "No active workbook is available."
    }

    # The template lives on the sheet named "template"; fall back to the
    # active sheet if, for some reason, that sheet cannot be found.
    $ws = $null
    foreach ($sheet in $wb.Worksheets) {
        if ($sheet.Name -eq "template") {
            $ws = $sheet
            break
        }
    }
    if (-not $ws) {
        $ws = $wb.ActiveSheet
    }
    if (-not $ws) {
        throw "Could not resolve a worksheet to edit."
    }

    # Remove the one-off "Påskmarknad 2022" title from B2. Keep the
    # cell's style/formatting intact; only the text content is cleared,
    # generalizing the template for reuse across events.
    $titleCell = $ws.Range("B2")
    $titleCell.ClearContents()

    # Reset the sheet's active selection to B2 (the now-blank title
    # cell) so the template opens with a sensible default selection.
    $ws.Activate()
    $titleCell.Select()
}
catch {
    Write-Host "Failed to generalize hämtlista template: $($_.Exception.Message)"
    throw
}
